$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 8497.75
$ws.Range("J51").Value = 8497.75
$ws.Range("L51").Value = 8497.75
$ws.Range("N51").Value = -9465.75
$ws.Range("H69").Value = 47625624
$ws.Range("J69").Value = 7669.1665
$ws.Range("L69").Value = 23007.4995
$ws.Range("N69").Value = -24755.4995
$ws.Range("H72").Value = 47625624
$ws.Range("J72").Value = 7669.1665
$ws.Range("L72").Value = 69022.4985
$ws.Range("N72").Value = -77758.4985
$ws.Range("H111").Value = 5996.3184
$ws.Range("I111").Value = 4715.6924
$ws.Range("J111").Value = 7846.1113
$ws.Range("K111").Value = 14147.0772
$ws.Range("L111").Value = 23538.3339
$ws.Range("M111").Value = -11080.0772
$ws.Range("N111").Value = -29672.3339
$ws.Range("H112").Value = 2530.8071
$ws.Range("J112").Value = 2572.2727
$ws.Range("L112").Value = 7716.8181
$ws.Range("N112").Value = -9932.8181
$ws.Range("H132").Value = 3401.8838
$ws.Range("I132").Value = 1807.175
$ws.Range("K132").Value = 5421.525
$ws.Range("M132").Value = -2891.525
$ws.Range("H133").Value = 116359.91
$ws.Range("J133").Value = 116359.91
$ws.Range("L133").Value = 116359.91
$ws.Range("N133").Value = -126479.91
$ws.Range("H137").Value = 15876164
$ws.Range("I137").Value = 1998.8572
$ws.Range("K137").Value = 5996.571599999999
$ws.Range("M137").Value = -3446.571599999999
$ws.Range("H141").Value = 3480
$ws.Range("I141").Value = 2647.8333
$ws.Range("K141").Value = 7943.499899999999
$ws.Range("M141").Value = -2763.499899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8508.521000000001
$ws.Range("I32").Value = 4527.2183
$ws.Range("K32").Value = 4527.2183
$ws.Range("M32").Value = -4240.2183
$ws.Range("H61").Value = 6277.6
$ws.Range("I61").Value = 6188.6523
$ws.Range("K61").Value = 6188.6523
$ws.Range("M61").Value = -5976.6523
$ws.Range("H74").Value = 31252874
$ws.Range("I74").Value = 62501504
$ws.Range("J74").Value = 4246
$ws.Range("K74").Value = 62501504
$ws.Range("L74").Value = 4246
$ws.Range("M74").Value = -62500630
$ws.Range("N74").Value = -5994
$ws.Range("H77").Value = 31252874
$ws.Range("I77").Value = 62501504
$ws.Range("J77").Value = 4246
$ws.Range("K77").Value = 312507520
$ws.Range("L77").Value = 21230
$ws.Range("M77").Value = -312503152
$ws.Range("N77").Value = -29966
$ws.Range("H110").Value = 2315.7693
$ws.Range("I110").Value = 2118.818
$ws.Range("J110").Value = 3399
$ws.Range("K110").Value = 2118.818
$ws.Range("L110").Value = 3399
$ws.Range("M110").Value = -73.81800000000021
$ws.Range("N110").Value = -7489
$ws.Range("H132").Value = 13209.451
$ws.Range("I132").Value = 16905.79
$ws.Range("J132").Value = 7356.9165
$ws.Range("K132").Value = 50717.37
$ws.Range("L132").Value = 22070.7495
$ws.Range("M132").Value = -48187.37
$ws.Range("N132").Value = -27130.7495
$ws.Range("H136").Value = 6277.6
$ws.Range("I136").Value = 6188.6523
$ws.Range("K136").Value = 18565.9569
$ws.Range("M136").Value = -16015.9569

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 21126.5
$ws.Range("I82").Value = 13858.857
$ws.Range("K82").Value = 13858.857
$ws.Range("M82").Value = -13475.857
$ws.Range("H85").Value = 21126.5
$ws.Range("I85").Value = 13858.857
$ws.Range("K85").Value = 13858.857
$ws.Range("M85").Value = -12532.857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17860162
$ws.Range("I31").Value = 26317808
$ws.Range("K31").Value = 26317808
$ws.Range("M31").Value = -26317513
$ws.Range("H34").Value = 17860162
$ws.Range("I34").Value = 26317808
$ws.Range("K34").Value = 26317808
$ws.Range("M34").Value = -26317606
$ws.Range("H58").Value = 3312.7693
$ws.Range("I58").Value = 2822.7778
$ws.Range("J58").Value = 4415.25
$ws.Range("K58").Value = 2822.7778
$ws.Range("L58").Value = 4415.25
$ws.Range("M58").Value = -2619.7778
$ws.Range("N58").Value = -4821.25
$ws.Range("H99").Value = 5466.846
$ws.Range("I99").Value = 2899.2856
$ws.Range("K99").Value = 2899.2856
$ws.Range("M99").Value = -1401.2856
$ws.Range("H122").Value = 3455.121
$ws.Range("I122").Value = 2123.3684
$ws.Range("J122").Value = 5262.5
$ws.Range("K122").Value = 6370.1052
$ws.Range("L122").Value = 15787.5
$ws.Range("M122").Value = -3920.1052
$ws.Range("N122").Value = -20687.5
$ws.Range("H126").Value = 5466.846
$ws.Range("I126").Value = 2899.2856
$ws.Range("K126").Value = 8697.856800000001
$ws.Range("M126").Value = -6227.856800000001
$ws.Range("H134").Value = 3988.7
$ws.Range("J134").Value = 4616.75
$ws.Range("L134").Value = 13850.25
$ws.Range("N134").Value = -18920.25
$ws.Range("H136").Value = 3312.7693
$ws.Range("I136").Value = 2822.7778
$ws.Range("J136").Value = 4415.25
$ws.Range("K136").Value = 8468.3334
$ws.Range("L136").Value = 13245.75
$ws.Range("M136").Value = -5918.3334
$ws.Range("N136").Value = -18345.75
$ws.Range("H141").Value = 307160.16
$ws.Range("J141").Value = 321957.34
$ws.Range("L141").Value = 321957.34
$ws.Range("N141").Value = -332317.34

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 333426660
$ws.Range("J37").Value = 333426660
$ws.Range("L37").Value = 1000279980
$ws.Range("N37").Value = -1000280204
$ws.Range("H122").Value = 1599.5
$ws.Range("I122").Value = 799.5
$ws.Range("J122").Value = 2399.5
$ws.Range("K122").Value = 7195.5
$ws.Range("L122").Value = 21595.5
$ws.Range("M122").Value = -4745.5
$ws.Range("N122").Value = -26495.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 945.9
$ws.Range("I107").Value = 451
$ws.Range("J107").Value = 1688.25
$ws.Range("K107").Value = 451
$ws.Range("L107").Value = 1688.25
$ws.Range("M107").Value = 1469
$ws.Range("N107").Value = -5528.25
$ws.Range("H132").Value = 3533
$ws.Range("I132").Value = 2983.12
$ws.Range("J132").Value = 5060.4443
$ws.Range("K132").Value = 8949.360000000001
$ws.Range("L132").Value = 15181.3329
$ws.Range("M132").Value = -6419.360000000001
$ws.Range("N132").Value = -20241.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 891.5454999999999
$ws.Range("I22").Value = 768.4706
$ws.Range("J22").Value = 1310
$ws.Range("K22").Value = 768.4706
$ws.Range("L22").Value = 1310
$ws.Range("M22").Value = -473.4706
$ws.Range("N22").Value = -1900
$ws.Range("H27").Value = 891.5454999999999
$ws.Range("I27").Value = 768.4706
$ws.Range("J27").Value = 1310
$ws.Range("K27").Value = 768.4706
$ws.Range("L27").Value = 1310
$ws.Range("M27").Value = -661.4706
$ws.Range("N27").Value = -1524
$ws.Range("H102").Value = 119748.75
$ws.Range("J102").Value = 119748.75
$ws.Range("L102").Value = 119748.75
$ws.Range("N102").Value = -126238.75
$ws.Range("H132").Value = 5731.863
$ws.Range("I132").Value = 5032.0815
$ws.Range("K132").Value = 15096.2445
$ws.Range("M132").Value = -12566.2445
$ws.Range("H134").Value = 94999.5
$ws.Range("J134").Value = 94999.5
$ws.Range("L134").Value = 94999.5
$ws.Range("N134").Value = -105139.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 50000
$ws.Range("I28").Value = 50000
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 50000
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -49652
$ws.Range("H132").Value = 3048.0212
$ws.Range("I132").Value = 1720.4546
$ws.Range("J132").Value = 4216.28
$ws.Range("K132").Value = 5161.3638
$ws.Range("L132").Value = 12648.84
$ws.Range("M132").Value = -2631.3638
$ws.Range("N132").Value = -17708.84
$ws.Range("H136").Value = 4800.8076
$ws.Range("I136").Value = 2797
$ws.Range("K136").Value = 8391
$ws.Range("M136").Value = -5841
$ws.Range("N28").ClearContents()
